# Tune Random forest algorithm
# Update predicted values in column B (Predicted Eg) for the rows
# whose predictions changed after retuning the model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 3.11
    3  = 3.31
    4  = 3.23
    5  = 3.24
    7  = 3.27
    8  = 3.19
    9  = 3.04
    10 = 3.23
    11 = 3.14
    12 = 3.15
    13 = 3.11
    14 = 3.1
    15 = 3.16
    16 = 3.13
    17 = 3.08
    18 = 3.04
    19 = 3.07
    20 = 1.92
    21 = 2.35
    22 = 2.27
    23 = 2.28
    24 = 2.25
    25 = 2.32
    26 = 6.87
    27 = 3.56
    28 = 3.78
    29 = 2.62
    31 = 8.5
    32 = 2.62
    33 = 2.66
    34 = 3
    35 = 3.8
    38 = 3.45
    39 = 3.47
    40 = 3.45
    41 = 3.44
    42 = 2.73
    43 = 2.7
    44 = 3.39
    46 = 3.37
    47 = 3.35
    48 = 2.62
    49 = 2.62
    50 = 2.66
    51 = 2.71
    52 = 2.68
    53 = 3.15
    54 = 3.17
    55 = 3.19
    56 = 3.24
    57 = 3.2
    58 = 3.23
    59 = 3.21
    60 = 3.2
    61 = 3.2
    62 = 3.2
    63 = 3.21
    64 = 3.21
    65 = 3.3
    66 = 3.22
    67 = 3.23
    68 = 3.19
    69 = 3.2
    70 = 3.2
    71 = 3.17
    73 = 3.3
    74 = 3.27
    75 = 3.25
    76 = 3.45
    77 = 3.41
    78 = 3.22
    79 = 3.27
    80 = 3.23
    81 = 3.22
    82 = 3.21
    85 = 3.12
    86 = 3.12
    87 = 3.31
    88 = 3.27
    89 = 3.22
    90 = 3.19
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}
